$d = $word.ActiveDocument

function New-PkgXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParaStartingWith($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($text)) { return $d.Paragraphs.Item($i) }
    }
    return $null
}

$ilvl1pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# ---------------------------------------------------------------------------
# Step 1 (process bottom-up so earlier paragraph indices stay valid):
# Insert a brand-new bullet right after the "Logistic analysis" heading
# paragraph, without touching that paragraph at all:
#   "The long logit has a pseudo-r2 of .45, noticeably less than the weak
#   linear model's r2."
# ---------------------------------------------------------------------------
$p25 = Find-ParaStartingWith("Logistic analysis")
$r25 = $p25.Range
$insAfter25 = $d.Range($r25.End, $r25.End)
$body1 = '<w:p>' + $ilvl1pPr + '<w:r><w:t>The long logit has a pseudo-r2 of .45, noticeably less than the weak linear model’s r2.</w:t></w:r></w:p>'
$insAfter25.InsertXML((New-PkgXml $body1))

# ---------------------------------------------------------------------------
# Step 2: The "Surprisingly, ndet..." paragraph loses its trailing
# bookmarkStart/bookmarkEnd for "_GoBack", and a new bullet is added right
# after it: "Q-complexity of the strong model is no different than the
# maxar model". Done as one full-paragraph replacement so the bookmark pair
# is cleanly dropped.
# ---------------------------------------------------------------------------
$p24 = Find-ParaStartingWith("Surprisingly, ndet")
$r24 = $p24.Range
$ndetText = "Surprisingly, ndet and lastupdate had strong effects. I interpret ndet as reflecting transparency and portfolio diversity. Factual transparency might be due to an individual’s desire to signal that they have much to offer, but transparency might also be attitudinal. Attitudinal transparency may be detected by an individual’s willingness to expose their portfolio even when the portfolio is unimpressive. Attitudinal transparency may be more a matter of personality traits than technical skill."
$body2 = '<w:p>' + $ilvl1pPr + '<w:r><w:t>' + $ndetText + '</w:t></w:r></w:p>' + `
         '<w:p>' + $ilvl1pPr + '<w:r><w:t>Q-complexity of the strong model is no different than the maxar model</w:t></w:r></w:p>'
$r24.InsertXML((New-PkgXml $body2))

# ---------------------------------------------------------------------------
# Step 3: The "Country effects weren't significant..." paragraph is replaced
# by two paragraphs:
#   (a) "Simple regression on country had an r2 of about .6, but country
#       effects were entirely omitted from non-simple standard reduction due
#       to collinearity. 78% of sampled users were from the US." plus the
#       _GoBack bookmark and a trailing space run.
#   (b) "Country effects " / "were omitted for collinearity. This may be
#       disconcerting" / ", possibly due to sampling, although sampling
#       well-represents the actual platform userbase. State effects were
#       significant, but only two states survived until the strong factor
#       model: CA and MI. These states had vertically-robust negative
#       effects."
# ---------------------------------------------------------------------------
$p21 = Find-ParaStartingWith("Country effects weren")
$r21 = $p21.Range
$tailText = ", possibly due to sampling, although sampling well-represents the actual platform userbase. State effects were significant, but only two states survived until the strong factor model: CA and MI. These states had vertically-robust negative effects."
$body3 = '<w:p>' + $ilvl1pPr + '<w:r><w:t>Simple regression on country had an r2 of about .6, but country effects were entirely omitted from non-simple standard reduction due to collinearity. 78% of sampled users were from the US.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + `
         '<w:p>' + $ilvl1pPr + '<w:r><w:t xml:space="preserve">Country effects </w:t></w:r><w:r><w:t>were omitted for collinearity. This may be disconcerting</w:t></w:r><w:r><w:t>' + $tailText + '</w:t></w:r></w:p>'
$r21.InsertXML((New-PkgXml $body3))

# ---------------------------------------------------------------------------
# Step 4: The "Structural effects..." paragraph gains two more runs:
#   " " and "Exploratory2-4 show structural importance of states,
#   interacted3, and nnano3."
# Replace only the text sub-range (not the <w:p> itself), so the
# paragraph's own opening tag is left completely untouched.
# ---------------------------------------------------------------------------
$p19 = Find-ParaStartingWith("Structural effects were lost")
$r19 = $p19.Range
$sub19 = $d.Range($r19.Start, $r19.End - 1)
$body4 = '<w:p><w:r><w:t>Structural effects were lost in the strong factor model, so the adjusted r2-maximizing model is generally preferred.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Exploratory2-4 show structural importance of states, interacted3, and nnano3.</w:t></w:r></w:p>'
$sub19.InsertXML((New-PkgXml $body4))

Write-Output "done"
